$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.220.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.747.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.737.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.50'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.482'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000253'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.378.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.749.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.347.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '498.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.720'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000134'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("E31").Value = '  +3.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.113'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.347'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.99%  '
$ws.Range("E39").Value = '  +4.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '459.07'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.943.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0359'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.37%  '
